$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date: 2016-09-02 02:19:30 -> 2016-09-02 02:20:24
$wsOverview.Range("G2").Value = "2016-09-02 02:20:24"
$wsOverview.Range("G4").Value = "2016-09-02 02:20:24"

# Priority ht -> mt (shared string also used by zh-cn and de-de sheets)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime 2016-09-02 02:19:25 -> 2016-09-02 02:20:19
$wsZhCn.Range("H2").Value = "2016-09-02 02:20:19"
$wsZhCn.Range("H4").Value = "2016-09-02 02:20:19"

# zh-cn sheet: Correspond Handback DateTime 2016-09-02 02:19:42 -> 2016-09-02 02:20:38
$wsZhCn.Range("K2").Value = "2016-09-02 02:20:38"
$wsZhCn.Range("K4").Value = "2016-09-02 02:20:38"

# de-de sheet: Correspond Handoff Datetime 2016-09-02 02:19:30 -> 2016-09-02 02:20:24 (shared with Overview)
$wsDeDe.Range("H2").Value = "2016-09-02 02:20:24"
$wsDeDe.Range("H4").Value = "2016-09-02 02:20:24"

# de-de sheet: Correspond Handback DateTime 2016-09-02 02:19:49 -> 2016-09-02 02:20:45
$wsDeDe.Range("K2").Value = "2016-09-02 02:20:45"
$wsDeDe.Range("K4").Value = "2016-09-02 02:20:45"
